$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: strip the green/theme-5 fill + centering from A2:C2 (back to plain default) ---
$ws.Range("A2:C2").ClearFormats()

# --- Row 3: B3 keeps its text ("id2") - no change needed ---

# --- Row 4: A4 keeps "key"; B4 text changes from "INT" to "key" ---
$ws.Range("B4").Value = "key"

# --- Rows 5-8: renumber the A-column counters (0,1,2,3 -> 1,2,3,4) ---
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 4

# --- Row 7: E7 text changes and gets wrap + taller row ---
$ws.Range("E7").Value = "[1,3,`n4,6]"
$ws.Range("E7").VerticalAlignment = -4108
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").WrapText = $true
$ws.Rows(7).RowHeight = 28.5

# --- Row 8: E8 text changes and gets wrap + taller row ---
$ws.Range("E8").Value = "[{1,3},`n{2,4}]"
$ws.Range("E8").VerticalAlignment = -4108
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E8").WrapText = $true
$ws.Rows(8).RowHeight = 28.5

# --- Row 11: drop the stray B11 cell, A11 keeps "KR" ---
$ws.Range("B11").Clear()

# --- Row 15: drop the stray B15 cell, A15 keeps "TW" ---
$ws.Range("B15").Clear()

# --- Rows 16-17: renumber A16 (0 -> 1) and insert a new row 17 ---
$ws.Range("A16").Value = 1
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = 12
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 302
$ws.Range("E17").Value = "[{1,2}]"
$ws.Range("A17:E17").VerticalAlignment = -4108
$ws.Range("A17:E17").HorizontalAlignment = -4108

# --- Row 19: extend span to A:E (cosmetic; content unchanged: "Th") ---

# --- Selection moves from B5 to B13 ---
$ws.Range("B13").Select()

Write-Output "done"
